# Regenerate handoff report: rename the source file's generated UUID / content hash,
# refresh handoff timestamps, and drop the row describing the previous failed
# handoff transform (which no longer applies) from every sheet.

$wb = $excel.ActiveWorkbook

$oldUuid = "7199f349-841a-4fa6-a12d-e45bc2f92b05"
$newUuid = "a0c99a9f-fea3-4945-bb44-704596b416ed"
$oldHash = "a6e78383983ac1dce43127ec7bf243e76975b61e"
$newHash = "b1d049f8c9afb7c5d5ab1424b52f66de5f5015b8"

$oldZhDate = "2016-01-14 15:54:57"
$newZhDate = "2016-01-14 15:57:12"
$oldDeDate = "2016-01-14 15:55:39"
$newDeDate = "2016-01-14 15:58:14"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

foreach ($ws in $wb.Worksheets) {
    $null = $ws.Cells.Replace($oldUuid, $newUuid)
    $null = $ws.Cells.Replace($oldHash, $newHash)
}
$null = $wsZh.Cells.Replace($oldZhDate, $newZhDate)
$null = $wsDe.Cells.Replace($oldDeDate, $newDeDate)

# The row describing "1b30fc14-...md" / "Handoff transform failed" is gone in the
# new report - remove row 3 (shifting the ".localization-config" row up to row 3)
# on every sheet.
$wsOverview.Rows("3").Delete()
$wsZh.Rows("3").Delete()
$wsDe.Rows("3").Delete()

function Clear-AllHyperlinks($ws) {
    while ($ws.Hyperlinks.Count -gt 0) {
        foreach ($h in $ws.Hyperlinks) {
            $h.Delete()
        }
    }
}

Clear-AllHyperlinks $wsOverview
Clear-AllHyperlinks $wsZh
Clear-AllHyperlinks $wsDe

$baseRepo = "https://github.com/OpenLocalizationTest/oltest/blob/ccdf518dbfad94bd0d4dffa363b40b7e3c6a18ad"
$mdUrl = $baseRepo + "/e2e/" + $newUuid + ".md"
$configUrl = $baseRepo + "/.localization-config"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa877ef9c9430120da42ae16c0337c378c46bc6d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/" + $newUuid + "." + $newHash + ".zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cdf843abc07cc93a19321c6cd8ba30b4e99a7815/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/" + $newUuid + "." + $newHash + ".de-de.xlf"

# Overview sheet: A2 -> md file, A3 -> .localization-config
$null = $wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $newUuid + ".md")
$null = $wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", ".localization-config")

# zh-cn sheet: A2 -> md file, C2 -> xlf handoff file, A3 -> .localization-config
$null = $wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", $newUuid + ".md")
$null = $wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfUrl, "", "", $newUuid + "." + $newHash + ".zh-cn.xlf")
$null = $wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configUrl, "", "", ".localization-config")

# de-de sheet: A2 -> md file, C2 -> xlf handoff file, A3 -> .localization-config
$null = $wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", $newUuid + ".md")
$null = $wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfUrl, "", "", $newUuid + "." + $newHash + ".de-de.xlf")
$null = $wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configUrl, "", "", ".localization-config")
